# Daily attendance processing - 2025-11-26 04:33:06
# Rotate the "Recorded By" (column G) comma-separated list of recorders
# left by one position (move the first name to the end of the list) for
# every data row in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value()

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
